{"js": "// Replace each old equation string with its new value.\n// The document body is a table of 100 unique \"a op b = c\" strings;\n// each (before, after) pair below is unique and unambiguous, so a\n// plain text search-and-replace keyed on the old value is safe.\nconst replacements = [\n  [\"65-19=46\", \"71-53=18\"],\n  [\"59-17=42\", \"85-75=10\"],\n  [\"31-20=11\", \"19-9=10\"],\n  [\"67-47=20\", \"37-17=20\"],\n  [\"5+34=39\", \"99-37=62\"],\n  [\"46+53=99\", \"70-33=37\"],\n  [\"73-64=9\", \"45-28=17\"],\n  [\"79-54=25\", \"78-43=35\"],\n  [\"28-18=10\", \"56-29=27\"],\n  [\"75+8=83\", \"33+54=87\"],\n  [\"33-3=30\", \"51-2=49\"],\n  [\"30-4=26\", \"1+32=33\"],\n  [\"96-54=42\", \"99-93=6\"],\n  [\"19+2=21\", \"1+40=41\"],\n  [\"16+8=24\", \"29+62=91\"],\n  [\"35-16=19\", \"47-30=17\"],\n  [\"90-81=9\", \"6+43=49\"],\n  [\"93-20=73\", \"99-68=31\"],\n  [\"33-19=14\", \"43-2=41\"],\n  [\"4+78=82\", \"91-16=75\"],\n  [\"18+26=44\", \"93-2=91\"],\n  [\"31-7=24\", \"23+6=29\"],\n  [\"78-12=66\", \"9-0=9\"],\n  [\"60-14=46\", \"89-1=88\"],\n  [\"16+37=53\", \"98-45=53\"],\n  [\"23+39=62\", \"59-23=36\"],\n  [\"53+5=58\", \"17+31=48\"],\n  [\"90-4=86\", \"97-31=66\"],\n  [\"69+19=88\", \"61+10=71\"],\n  [\"59-39=20\", \"4+64=68\"],\n  [\"18+39=57\", \"76-55=21\"],\n  [\"10+83=93\", \"89-16=73\"],\n  [\"91-80=11\", \"62-42=20\"],\n  [\"2+74=76\", \"43-20=23\"],\n  [\"1+5=6\", \"51+18=69\"],\n  [\"43-39=4\", \"42-1=41\"],\n  [\"59-5=54\", \"91-38=53\"],\n  [\"26+60=86\", \"6+70=76\"],\n  [\"31+50=81\", \"71-35=36\"],\n  [\"69-26=43\", \"66-34=32\"],\n  [\"24-12=12\", \"8+1=9\"],\n  [\"17+29=46\", \"61-38=23\"],\n  [\"77-15=62\", \"47-16=31\"],\n  [\"0+70=70\", \"15+36=51\"],\n  [\"96-27=69\", \"66-43=23\"],\n  [\"58+22=80\", \"5-1=4\"],\n  [\"60+33=93\", \"86+13=99\"],\n  [\"96-69=27\", \"73-6=67\"],\n  [\"48+48=96\", \"13+55=68\"],\n  [\"78+18=96\", \"90-7=83\"],\n  [\"86-28=58\", \"38+0=38\"],\n  [\"83-2=81\", \"14+44=58\"],\n  [\"25+10=35\", \"77-25=52\"],\n  [\"47+41=88\", \"84+6=90\"],\n  [\"30+25=55\", \"40-14=26\"],\n  [\"53-36=17\", \"48+34=82\"],\n  [\"78-45=33\", \"25+26=51\"],\n  [\"44+33=77\", \"71+8=79\"],\n  [\"93-23=70\", \"6+9=15\"],\n  [\"65+6=71\", \"15+43=58\"],\n  [\"93-74=19\", \"1+0=1\"],\n  [\"88-14=74\", \"18-16=2\"],\n  [\"94-12=82\", \"50-18=32\"],\n  [\"35+41=76\", \"24+25=49\"],\n  [\"40+34=74\", \"72-45=27\"],\n  [\"11+6=17\", \"4+30=34\"],\n  [\"8+84=92\", \"85-71=14\"],\n  [\"1+38=39\", \"75-75=0\"],\n  [\"12-8=4\", \"64-63=1\"],\n  [\"39-15=24\", \"57-20=37\"],\n  [\"17+6=23\", \"56+26=82\"],\n  [\"92-73=19\", \"25-10=15\"],\n  [\"43+4=47\", \"21+21=42\"],\n  [\"41+58=99\", \"63+10=73\"],\n  [\"88-32=56\", \"55-26=29\"],\n  [\"48+20=68\", \"54-32=22\"],\n  [\"6+38=44\", \"23+31=54\"],\n  [\"71-5=66\", \"33+55=88\"],\n  [\"35+9=44\", \"93-52=41\"],\n  [\"79-44=35\", \"88-64=24\"],\n  [\"81-43=38\", \"95-36=59\"],\n  [\"0+64=64\", \"48+34=82\"],\n  [\"29+27=56\", \"93-71=22\"],\n  [\"68-19=49\", \"19+9=28\"],\n  [\"4-3=1\", \"39-36=3\"],\n  [\"30-24=6\", \"40-18=22\"],\n  [\"37+15=52\", \"16-5=11\"],\n  [\"94-42=52\", \"35-32=3\"],\n  [\"80-57=23\", \"20+13=33\"],\n  [\"51-16=35\", \"61+28=89\"],\n  [\"31+13=44\", \"17+72=89\"],\n  [\"33+28=61\", \"45+9=54\"],\n  [\"11+38=49\", \"95-28=67\"],\n  [\"48-17=31\", \"55-46=9\"],\n  [\"88-41=47\", \"36+1=37\"],\n  [\"99-79=20\", \"82-67=15\"],\n  [\"55-53=2\", \"12+42=54\"],\n  [\"30-5=25\", \"46+37=83\"],\n  [\"13+53=66\", \"78-60=18\"],\n  [\"28+28=56\", \"83-59=24\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Pattern not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each old equation string with its new value.\n# The document body is a table of 100 unique \"a op b = c\" strings;\n# each (Find, Replace) pair below is unique and unambiguous, so a\n# plain text Find/Replace keyed on the old value is safe and\n# unambiguous regardless of iteration order.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"65-19=46\"; Replace = \"71-53=18\" },\n    @{ Find = \"59-17=42\"; Replace = \"85-75=10\" },\n    @{ Find = \"31-20=11\"; Replace = \"19-9=10\" },\n    @{ Find = \"67-47=20\"; Replace = \"37-17=20\" },\n    @{ Find = \"5+34=39\"; Replace = \"99-37=62\" },\n    @{ Find = \"46+53=99\"; Replace = \"70-33=37\" },\n    @{ Find = \"73-64=9\"; Replace = \"45-28=17\" },\n    @{ Find = \"79-54=25\"; Replace = \"78-43=35\" },\n    @{ Find = \"28-18=10\"; Replace = \"56-29=27\" },\n    @{ Find = \"75+8=83\"; Replace = \"33+54=87\" },\n    @{ Find = \"33-3=30\"; Replace = \"51-2=49\" },\n    @{ Find = \"30-4=26\"; Replace = \"1+32=33\" },\n    @{ Find = \"96-54=42\"; Replace = \"99-93=6\" },\n    @{ Find = \"19+2=21\"; Replace = \"1+40=41\" },\n    @{ Find = \"16+8=24\"; Replace = \"29+62=91\" },\n    @{ Find = \"35-16=19\"; Replace = \"47-30=17\" },\n    @{ Find = \"90-81=9\"; Replace = \"6+43=49\" },\n    @{ Find = \"93-20=73\"; Replace = \"99-68=31\" },\n    @{ Find = \"33-19=14\"; Replace = \"43-2=41\" },\n    @{ Find = \"4+78=82\"; Replace = \"91-16=75\" },\n    @{ Find = \"18+26=44\"; Replace = \"93-2=91\" },\n    @{ Find = \"31-7=24\"; Replace = \"23+6=29\" },\n    @{ Find = \"78-12=66\"; Replace = \"9-0=9\" },\n    @{ Find = \"60-14=46\"; Replace = \"89-1=88\" },\n    @{ Find = \"16+37=53\"; Replace = \"98-45=53\" },\n    @{ Find = \"23+39=62\"; Replace = \"59-23=36\" },\n    @{ Find = \"53+5=58\"; Replace = \"17+31=48\" },\n    @{ Find = \"90-4=86\"; Replace = \"97-31=66\" },\n    @{ Find = \"69+19=88\"; Replace = \"61+10=71\" },\n    @{ Find = \"59-39=20\"; Replace = \"4+64=68\" },\n    @{ Find = \"18+39=57\"; Replace = \"76-55=21\" },\n    @{ Find = \"10+83=93\"; Replace = \"89-16=73\" },\n    @{ Find = \"91-80=11\"; Replace = \"62-42=20\" },\n    @{ Find = \"2+74=76\"; Replace = \"43-20=23\" },\n    @{ Find = \"1+5=6\"; Replace = \"51+18=69\" },\n    @{ Find = \"43-39=4\"; Replace = \"42-1=41\" },\n    @{ Find = \"59-5=54\"; Replace = \"91-38=53\" },\n    @{ Find = \"26+60=86\"; Replace = \"6+70=76\" },\n    @{ Find = \"31+50=81\"; Replace = \"71-35=36\" },\n    @{ Find = \"69-26=43\"; Replace = \"66-34=32\" },\n    @{ Find = \"24-12=12\"; Replace = \"8+1=9\" },\n    @{ Find = \"17+29=46\"; Replace = \"61-38=23\" },\n    @{ Find = \"77-15=62\"; Replace = \"47-16=31\" },\n    @{ Find = \"0+70=70\"; Replace = \"15+36=51\" },\n    @{ Find = \"96-27=69\"; Replace = \"66-43=23\" },\n    @{ Find = \"58+22=80\"; Replace = \"5-1=4\" },\n    @{ Find = \"60+33=93\"; Replace = \"86+13=99\" },\n    @{ Find = \"96-69=27\"; Replace = \"73-6=67\" },\n    @{ Find = \"48+48=96\"; Replace = \"13+55=68\" },\n    @{ Find = \"78+18=96\"; Replace = \"90-7=83\" },\n    @{ Find = \"86-28=58\"; Replace = \"38+0=38\" },\n    @{ Find = \"83-2=81\"; Replace = \"14+44=58\" },\n    @{ Find = \"25+10=35\"; Replace = \"77-25=52\" },\n    @{ Find = \"47+41=88\"; Replace = \"84+6=90\" },\n    @{ Find = \"30+25=55\"; Replace = \"40-14=26\" },\n    @{ Find = \"53-36=17\"; Replace = \"48+34=82\" },\n    @{ Find = \"78-45=33\"; Replace = \"25+26=51\" },\n    @{ Find = \"44+33=77\"; Replace = \"71+8=79\" },\n    @{ Find = \"93-23=70\"; Replace = \"6+9=15\" },\n    @{ Find = \"65+6=71\"; Replace = \"15+43=58\" },\n    @{ Find = \"93-74=19\"; Replace = \"1+0=1\" },\n    @{ Find = \"88-14=74\"; Replace = \"18-16=2\" },\n    @{ Find = \"94-12=82\"; Replace = \"50-18=32\" },\n    @{ Find = \"35+41=76\"; Replace = \"24+25=49\" },\n    @{ Find = \"40+34=74\"; Replace = \"72-45=27\" },\n    @{ Find = \"11+6=17\"; Replace = \"4+30=34\" },\n    @{ Find = \"8+84=92\"; Replace = \"85-71=14\" },\n    @{ Find = \"1+38=39\"; Replace = \"75-75=0\" },\n    @{ Find = \"12-8=4\"; Replace = \"64-63=1\" },\n    @{ Find = \"39-15=24\"; Replace = \"57-20=37\" },\n    @{ Find = \"17+6=23\"; Replace = \"56+26=82\" },\n    @{ Find = \"92-73=19\"; Replace = \"25-10=15\" },\n    @{ Find = \"43+4=47\"; Replace = \"21+21=42\" },\n    @{ Find = \"41+58=99\"; Replace = \"63+10=73\" },\n    @{ Find = \"88-32=56\"; Replace = \"55-26=29\" },\n    @{ Find = \"48+20=68\"; Replace = \"54-32=22\" },\n    @{ Find = \"6+38=44\"; Replace = \"23+31=54\" },\n    @{ Find = \"71-5=66\"; Replace = \"33+55=88\" },\n    @{ Find = \"35+9=44\"; Replace = \"93-52=41\" },\n    @{ Find = \"79-44=35\"; Replace = \"88-64=24\" },\n    @{ Find = \"81-43=38\"; Replace = \"95-36=59\" },\n    @{ Find = \"0+64=64\"; Replace = \"48+34=82\" },\n    @{ Find = \"29+27=56\"; Replace = \"93-71=22\" },\n    @{ Find = \"68-19=49\"; Replace = \"19+9=28\" },\n    @{ Find = \"4-3=1\"; Replace = \"39-36=3\" },\n    @{ Find = \"30-24=6\"; Replace = \"40-18=22\" },\n    @{ Find = \"37+15=52\"; Replace = \"16-5=11\" },\n    @{ Find = \"94-42=52\"; Replace = \"35-32=3\" },\n    @{ Find = \"80-57=23\"; Replace = \"20+13=33\" },\n    @{ Find = \"51-16=35\"; Replace = \"61+28=89\" },\n    @{ Find = \"31+13=44\"; Replace = \"17+72=89\" },\n    @{ Find = \"33+28=61\"; Replace = \"45+9=54\" },\n    @{ Find = \"11+38=49\"; Replace = \"95-28=67\" },\n    @{ Find = \"48-17=31\"; Replace = \"55-46=9\" },\n    @{ Find = \"88-41=47\"; Replace = \"36+1=37\" },\n    @{ Find = \"99-79=20\"; Replace = \"82-67=15\" },\n    @{ Find = \"55-53=2\"; Replace = \"12+42=54\" },\n    @{ Find = \"30-5=25\"; Replace = \"46+37=83\" },\n    @{ Find = \"13+53=66\"; Replace = \"78-60=18\" },\n    @{ Find = \"28+28=56\"; Replace = \"83-59=24\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Find\n    $find.Replacement.Text = $pair.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    $found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n    if (-not $found) {\n        throw \"Pattern not found: $($pair.Find)\"\n    }\n}\n"}
